# Auto-generated Excel COM-interop script applying scheduled market-data refresh
# to the Anima_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
#
# For each changed row we update currentAveragePrice / currentAveragePriceNQ /
# currentAveragePriceHQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ
# (columns H-N) with refreshed market values. A couple of rows gain or lose a
# profit cell entirely (market data no longer / newly available for that recipe).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1321.83  # H15: 634.24 -> 1321.83
$ws.Cells.Item(15, 9).Value = 1321.83  # I15: 634.24 -> 1321.83
$ws.Cells.Item(15, 11).Value = 3965.49  # K15: 1902.72 -> 3965.49
$ws.Cells.Item(15, 13).Value = -3796.49  # M15: -1733.72 -> -3796.49

$ws.Cells.Item(22, 8).Value = 1243  # H22: 1500 -> 1243
$ws.Cells.Item(22, 9).Value = 472  # I22: 0 -> 472
$ws.Cells.Item(22, 11).Value = 1416  # K22: 0 -> 1416
$ws.Cells.Item(22, 13).Value = -1244  # M22: None -> -1244

$ws.Cells.Item(99, 8).Value = 2072.2222  # H99: 1820.375 -> 2072.2222
$ws.Cells.Item(99, 9).Value = 714.6667  # I99: 626 -> 714.6667
$ws.Cells.Item(99, 10).Value = 4787.3335  # J99: 10181 -> 4787.3335
$ws.Cells.Item(99, 11).Value = 2144.0001  # K99: 1878 -> 2144.0001
$ws.Cells.Item(99, 12).Value = 14362.0005  # L99: 30543 -> 14362.0005
$ws.Cells.Item(99, 13).Value = -646.0001000000002  # M99: -380 -> -646.0001000000002
$ws.Cells.Item(99, 14).Value = -17358.0005  # N99: -33539 -> -17358.0005

$ws.Cells.Item(116, 8).Value = 8706.066000000001  # H116: 9356.071 -> 8706.066000000001
$ws.Cells.Item(116, 9).Value = 12040.5  # I116: 13378.333 -> 12040.5
$ws.Cells.Item(116, 10).Value = 2037.2  # J116: 2116 -> 2037.2
$ws.Cells.Item(116, 11).Value = 12040.5  # K116: 13378.333 -> 12040.5
$ws.Cells.Item(116, 12).Value = 2037.2  # L116: 2116 -> 2037.2
$ws.Cells.Item(116, 13).Value = -8598.5  # M116: -9936.333000000001 -> -8598.5
$ws.Cells.Item(116, 14).Value = -8921.200000000001  # N116: -9000 -> -8921.200000000001

$ws.Cells.Item(132, 8).Value = 2716.5  # H132: 2988.0557 -> 2716.5
$ws.Cells.Item(132, 9).Value = 2622.5637  # I132: 2910.3264 -> 2622.5637
$ws.Cells.Item(132, 11).Value = 7867.6911  # K132: 8730.9792 -> 7867.6911
$ws.Cells.Item(132, 13).Value = -5337.6911  # M132: -6200.9792 -> -5337.6911

$ws.Cells.Item(135, 8).Value = 816.8461  # H135: 788.46295 -> 816.8461
$ws.Cells.Item(135, 9).Value = 453.27658  # I135: 436.83673 -> 453.27658
$ws.Cells.Item(135, 11).Value = 4079.48922  # K135: 3931.53057 -> 4079.48922
$ws.Cells.Item(135, 13).Value = -1544.48922  # M135: -1396.53057 -> -1544.48922

$ws.Cells.Item(138, 8).Value = 2098.2659  # H138: 2072.5212 -> 2098.2659
$ws.Cells.Item(138, 9).Value = 1399.7073  # I138: 1340.9773 -> 1399.7073
$ws.Cells.Item(138, 10).Value = 2638.6604  # J138: 2716.28 -> 2638.6604
$ws.Cells.Item(138, 11).Value = 4199.1219  # K138: 4022.9319 -> 4199.1219
$ws.Cells.Item(138, 12).Value = 7915.9812  # L138: 8148.84 -> 7915.9812
$ws.Cells.Item(138, 13).Value = 940.8780999999999  # M138: 1117.0681 -> 940.8780999999999
$ws.Cells.Item(138, 14).Value = -18195.9812  # N138: -18428.84 -> -18195.9812

$ws.Cells.Item(141, 8).Value = 4459.9  # H141: 4599.8965 -> 4459.9
$ws.Cells.Item(141, 10).Value = 9209.091  # J141: 10090 -> 9209.091
$ws.Cells.Item(141, 12).Value = 27627.273  # L141: 30270 -> 27627.273
$ws.Cells.Item(141, 14).Value = -37987.273  # N141: -40630 -> -37987.273

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(28, 8).Value = 3905.6  # H28: 3885.5833 -> 3905.6
$ws.Cells.Item(28, 9).Value = 3905.6  # I28: 3885.5833 -> 3905.6
$ws.Cells.Item(28, 11).Value = 3905.6  # K28: 3885.5833 -> 3905.6
$ws.Cells.Item(28, 13).Value = -3713.6  # M28: -3693.5833 -> -3713.6

$ws.Cells.Item(32, 8).Value = 650020.75  # H32: 692881.5600000001 -> 650020.75
$ws.Cells.Item(32, 9).Value = 739706.9  # I32: 795889.8 -> 739706.9
$ws.Cells.Item(32, 11).Value = 739706.9  # K32: 795889.8 -> 739706.9
$ws.Cells.Item(32, 13).Value = -739419.9  # M32: -795602.8 -> -739419.9

$ws.Cells.Item(61, 8).Value = 1789.5964  # H61: 2014.0869 -> 1789.5964
$ws.Cells.Item(61, 9).Value = 1248.0851  # I61: 1491.579 -> 1248.0851
$ws.Cells.Item(61, 10).Value = 4334.7  # J61: 4496 -> 4334.7
$ws.Cells.Item(61, 11).Value = 1248.0851  # K61: 1491.579 -> 1248.0851
$ws.Cells.Item(61, 12).Value = 4334.7  # L61: 4496 -> 4334.7
$ws.Cells.Item(61, 13).Value = -1036.0851  # M61: -1279.579 -> -1036.0851
$ws.Cells.Item(61, 14).Value = -4758.7  # N61: -4920 -> -4758.7

$ws.Cells.Item(74, 8).Value = 890.75757  # H74: 1226.1111 -> 890.75757
$ws.Cells.Item(74, 9).Value = 785.5357  # I74: 1009.75 -> 785.5357
$ws.Cells.Item(74, 10).Value = 1480  # J74: 2957 -> 1480
$ws.Cells.Item(74, 11).Value = 785.5357  # K74: 1009.75 -> 785.5357
$ws.Cells.Item(74, 12).Value = 1480  # L74: 2957 -> 1480
$ws.Cells.Item(74, 13).Value = 88.46429999999998  # M74: -135.75 -> 88.46429999999998
$ws.Cells.Item(74, 14).Value = -3228  # N74: -4705 -> -3228

$ws.Cells.Item(77, 8).Value = 890.75757  # H77: 1226.1111 -> 890.75757
$ws.Cells.Item(77, 9).Value = 785.5357  # I77: 1009.75 -> 785.5357
$ws.Cells.Item(77, 10).Value = 1480  # J77: 2957 -> 1480
$ws.Cells.Item(77, 11).Value = 3927.6785  # K77: 5048.75 -> 3927.6785
$ws.Cells.Item(77, 12).Value = 7400  # L77: 14785 -> 7400
$ws.Cells.Item(77, 13).Value = 440.3215  # M77: -680.75 -> 440.3215
$ws.Cells.Item(77, 14).Value = -16136  # N77: -23521 -> -16136

$ws.Cells.Item(99, 8).Value = 3905.6  # H99: 3885.5833 -> 3905.6
$ws.Cells.Item(99, 9).Value = 3905.6  # I99: 3885.5833 -> 3905.6
$ws.Cells.Item(99, 11).Value = 3905.6  # K99: 3885.5833 -> 3905.6
$ws.Cells.Item(99, 13).Value = -910.5999999999999  # M99: -890.5832999999998 -> -910.5999999999999

$ws.Cells.Item(113, 8).Value = 0  # H113: 44444 -> 0
$ws.Cells.Item(113, 10).Value = 0  # J113: 44444 -> 0
$ws.Cells.Item(113, 12).Value = 0  # L113: 44444 -> 0
$ws.Cells.Item(113, 14).ClearContents()  # N113: was -53122

$ws.Cells.Item(132, 8).Value = 3058.1836  # H132: 2226.7124 -> 3058.1836
$ws.Cells.Item(132, 9).Value = 2301.625  # I132: 1557.2909 -> 2301.625
$ws.Cells.Item(132, 10).Value = 4482.294  # J132: 4272.1665 -> 4482.294
$ws.Cells.Item(132, 11).Value = 6904.875  # K132: 4671.8727 -> 6904.875
$ws.Cells.Item(132, 12).Value = 13446.882  # L132: 12816.4995 -> 13446.882
$ws.Cells.Item(132, 13).Value = -4374.875  # M132: -2141.8727 -> -4374.875
$ws.Cells.Item(132, 14).Value = -18506.882  # N132: -17876.4995 -> -18506.882

$ws.Cells.Item(136, 8).Value = 1789.5964  # H136: 2014.0869 -> 1789.5964
$ws.Cells.Item(136, 9).Value = 1248.0851  # I136: 1491.579 -> 1248.0851
$ws.Cells.Item(136, 10).Value = 4334.7  # J136: 4496 -> 4334.7
$ws.Cells.Item(136, 11).Value = 3744.2553  # K136: 4474.737 -> 3744.2553
$ws.Cells.Item(136, 12).Value = 13004.1  # L136: 13488 -> 13004.1
$ws.Cells.Item(136, 13).Value = -1194.2553  # M136: -1924.737 -> -1194.2553
$ws.Cells.Item(136, 14).Value = -18104.1  # N136: -18588 -> -18104.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2041.2157  # H134: 2786.639 -> 2041.2157
$ws.Cells.Item(134, 9).Value = 1592  # I134: 2319.9614 -> 1592
$ws.Cells.Item(134, 10).Value = 3501.1667  # J134: 4000 -> 3501.1667
$ws.Cells.Item(134, 11).Value = 4776  # K134: 6959.8842 -> 4776
$ws.Cells.Item(134, 12).Value = 10503.5001  # L134: 12000 -> 10503.5001
$ws.Cells.Item(134, 13).Value = -2241  # M134: -4424.8842 -> -2241
$ws.Cells.Item(134, 14).Value = -15573.5001  # N134: -17070 -> -15573.5001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4664.0786  # H31: 5301.6045 -> 4664.0786
$ws.Cells.Item(31, 9).Value = 1444.1052  # I31: 1479.9445 -> 1444.1052
$ws.Cells.Item(31, 10).Value = 6575.9375  # J31: 8053.2 -> 6575.9375
$ws.Cells.Item(31, 11).Value = 1444.1052  # K31: 1479.9445 -> 1444.1052
$ws.Cells.Item(31, 12).Value = 6575.9375  # L31: 8053.2 -> 6575.9375
$ws.Cells.Item(31, 13).Value = -1149.1052  # M31: -1184.9445 -> -1149.1052
$ws.Cells.Item(31, 14).Value = -7165.9375  # N31: -8643.200000000001 -> -7165.9375

$ws.Cells.Item(34, 8).Value = 4664.0786  # H34: 5301.6045 -> 4664.0786
$ws.Cells.Item(34, 9).Value = 1444.1052  # I34: 1479.9445 -> 1444.1052
$ws.Cells.Item(34, 10).Value = 6575.9375  # J34: 8053.2 -> 6575.9375
$ws.Cells.Item(34, 11).Value = 1444.1052  # K34: 1479.9445 -> 1444.1052
$ws.Cells.Item(34, 12).Value = 6575.9375  # L34: 8053.2 -> 6575.9375
$ws.Cells.Item(34, 13).Value = -1242.1052  # M34: -1277.9445 -> -1242.1052
$ws.Cells.Item(34, 14).Value = -6979.9375  # N34: -8457.200000000001 -> -6979.9375

$ws.Cells.Item(58, 8).Value = 925.5893  # H58: 1056.7446 -> 925.5893
$ws.Cells.Item(58, 9).Value = 643.8378  # I58: 772.62067 -> 643.8378
$ws.Cells.Item(58, 10).Value = 1474.2632  # J58: 1514.5 -> 1474.2632
$ws.Cells.Item(58, 11).Value = 643.8378  # K58: 772.62067 -> 643.8378
$ws.Cells.Item(58, 12).Value = 1474.2632  # L58: 1514.5 -> 1474.2632
$ws.Cells.Item(58, 13).Value = -440.8378  # M58: -569.62067 -> -440.8378
$ws.Cells.Item(58, 14).Value = -1880.2632  # N58: -1920.5 -> -1880.2632

$ws.Cells.Item(118, 8).Value = 40000  # H118: 0 -> 40000
$ws.Cells.Item(118, 10).Value = 40000  # J118: 0 -> 40000
$ws.Cells.Item(118, 12).Value = 40000  # L118: 0 -> 40000
$ws.Cells.Item(118, 14).Value = -43314  # N118: None -> -43314

$ws.Cells.Item(132, 8).Value = 4387514.5  # H132: 3087642.8 -> 4387514.5
$ws.Cells.Item(132, 9).Value = 1474.68  # I132: 1073.5853 -> 1474.68
$ws.Cells.Item(132, 11).Value = 4424.04  # K132: 3220.7559 -> 4424.04
$ws.Cells.Item(132, 13).Value = -1894.04  # M132: -690.7559000000001 -> -1894.04

$ws.Cells.Item(134, 8).Value = 5676.7144  # H134: 3117.2632 -> 5676.7144
$ws.Cells.Item(134, 9).Value = 6397.15  # I134: 3047.0625 -> 6397.15
$ws.Cells.Item(134, 10).Value = 3875.625  # J134: 3491.6667 -> 3875.625
$ws.Cells.Item(134, 11).Value = 19191.45  # K134: 9141.1875 -> 19191.45
$ws.Cells.Item(134, 12).Value = 11626.875  # L134: 10475.0001 -> 11626.875
$ws.Cells.Item(134, 13).Value = -16656.45  # M134: -6606.1875 -> -16656.45
$ws.Cells.Item(134, 14).Value = -16696.875  # N134: -15545.0001 -> -16696.875

$ws.Cells.Item(136, 8).Value = 925.5893  # H136: 1056.7446 -> 925.5893
$ws.Cells.Item(136, 9).Value = 643.8378  # I136: 772.62067 -> 643.8378
$ws.Cells.Item(136, 10).Value = 1474.2632  # J136: 1514.5 -> 1474.2632
$ws.Cells.Item(136, 11).Value = 1931.5134  # K136: 2317.86201 -> 1931.5134
$ws.Cells.Item(136, 12).Value = 4422.7896  # L136: 4543.5 -> 4422.7896
$ws.Cells.Item(136, 13).Value = 618.4866  # M136: 232.1379900000002 -> 618.4866
$ws.Cells.Item(136, 14).Value = -9522.7896  # N136: -9643.5 -> -9522.7896

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 71428776  # H14: 71428824 -> 71428776
$ws.Cells.Item(14, 9).Value = 71428776  # I14: 71428824 -> 71428776
$ws.Cells.Item(14, 11).Value = 214286328  # K14: 214286472 -> 214286328
$ws.Cells.Item(14, 13).Value = -214286155  # M14: -214286299 -> -214286155

$ws.Cells.Item(17, 8).Value = 1749.5  # H17: 1249.6666 -> 1749.5
$ws.Cells.Item(17, 9).Value = 500  # I17: 375 -> 500
$ws.Cells.Item(17, 11).Value = 1500  # K17: 1125 -> 1500
$ws.Cells.Item(17, 13).Value = -1331  # M17: -956 -> -1331

$ws.Cells.Item(34, 8).Value = 11905239  # H34: 10417126 -> 11905239
$ws.Cells.Item(34, 9).Value = 211  # I34: 188.25 -> 211
$ws.Cells.Item(34, 10).Value = 12821010  # J34: 11364120 -> 12821010
$ws.Cells.Item(34, 11).Value = 633  # K34: 564.75 -> 633
$ws.Cells.Item(34, 12).Value = 38463030  # L34: 34092360 -> 38463030
$ws.Cells.Item(34, 13).Value = -549  # M34: -480.75 -> -549
$ws.Cells.Item(34, 14).Value = -38463198  # N34: -34092528 -> -38463198

$ws.Cells.Item(39, 8).Value = 1729.2916  # H39: 1700.1666 -> 1729.2916
$ws.Cells.Item(39, 9).Value = 1350  # I39: 1300 -> 1350
$ws.Cells.Item(39, 10).Value = 1763.7727  # J39: 1757.3334 -> 1763.7727
$ws.Cells.Item(39, 11).Value = 4050  # K39: 3900 -> 4050
$ws.Cells.Item(39, 12).Value = 5291.3181  # L39: 5272.0002 -> 5291.3181
$ws.Cells.Item(39, 13).Value = -3756  # M39: -3606 -> -3756
$ws.Cells.Item(39, 14).Value = -5879.3181  # N39: -5860.0002 -> -5879.3181

$ws.Cells.Item(55, 8).Value = 1505.8235  # H55: 1497 -> 1505.8235
$ws.Cells.Item(55, 10).Value = 1549.9375  # J55: 1540.5625 -> 1549.9375
$ws.Cells.Item(55, 12).Value = 4649.8125  # L55: 4621.6875 -> 4649.8125
$ws.Cells.Item(55, 14).Value = -5003.8125  # N55: -4975.6875 -> -5003.8125

$ws.Cells.Item(64, 8).Value = 2073.1428  # H64: 2064 -> 2073.1428
$ws.Cells.Item(64, 10).Value = 2083.3333  # J64: 2071.4285 -> 2083.3333
$ws.Cells.Item(64, 12).Value = 6249.999899999999  # L64: 6214.2855 -> 6249.999899999999
$ws.Cells.Item(64, 14).Value = -6789.999899999999  # N64: -6754.2855 -> -6789.999899999999

$ws.Cells.Item(67, 8).Value = 2073.1428  # H67: 2064 -> 2073.1428
$ws.Cells.Item(67, 10).Value = 2083.3333  # J67: 2071.4285 -> 2083.3333
$ws.Cells.Item(67, 12).Value = 6249.999899999999  # L67: 6214.2855 -> 6249.999899999999
$ws.Cells.Item(67, 14).Value = -8121.999899999999  # N67: -8086.2855 -> -8121.999899999999

$ws.Cells.Item(88, 8).Value = 7617.125  # H88: 7229.8184 -> 7617.125
$ws.Cells.Item(88, 10).Value = 7617.125  # J88: 7229.8184 -> 7617.125
$ws.Cells.Item(88, 12).Value = 22851.375  # L88: 21689.4552 -> 22851.375
$ws.Cells.Item(88, 14).Value = -23707.375  # N88: -22545.4552 -> -23707.375

$ws.Cells.Item(91, 8).Value = 7617.125  # H91: 7229.8184 -> 7617.125
$ws.Cells.Item(91, 10).Value = 7617.125  # J91: 7229.8184 -> 7617.125
$ws.Cells.Item(91, 12).Value = 22851.375  # L91: 21689.4552 -> 22851.375
$ws.Cells.Item(91, 14).Value = -25815.375  # N91: -24653.4552 -> -25815.375

$ws.Cells.Item(107, 8).Value = 20833602  # H107: 16666896 -> 20833602
$ws.Cells.Item(107, 10).Value = 62500236  # J107: 35714452 -> 62500236
$ws.Cells.Item(107, 12).Value = 187500708  # L107: 107143356 -> 187500708
$ws.Cells.Item(107, 14).Value = -187504548  # N107: -107147196 -> -187504548

$ws.Cells.Item(112, 8).Value = 4526.1665  # H112: 4286.3125 -> 4526.1665
$ws.Cells.Item(112, 9).Value = 3250.5715  # I112: 2642.3333 -> 3250.5715
$ws.Cells.Item(112, 10).Value = 6312  # J112: 6400 -> 6312
$ws.Cells.Item(112, 11).Value = 9751.7145  # K112: 7926.999899999999 -> 9751.7145
$ws.Cells.Item(112, 12).Value = 18936  # L112: 19200 -> 18936
$ws.Cells.Item(112, 13).Value = -8643.7145  # M112: -6818.999899999999 -> -8643.7145
$ws.Cells.Item(112, 14).Value = -21152  # N112: -21416 -> -21152

$ws.Cells.Item(132, 8).Value = 2001.3334  # H132: 1985.1154 -> 2001.3334
$ws.Cells.Item(132, 10).Value = 1869.2307  # J132: 1823.0834 -> 1869.2307
$ws.Cells.Item(132, 12).Value = 16823.0763  # L132: 16407.7506 -> 16823.0763
$ws.Cells.Item(132, 14).Value = -21883.0763  # N132: -21467.7506 -> -21883.0763

$ws.Cells.Item(139, 8).Value = 3368.0513  # H139: 2701.3845 -> 3368.0513
$ws.Cells.Item(139, 9).Value = 1762.3529  # I139: 1822.5 -> 1762.3529
$ws.Cells.Item(139, 10).Value = 4608.8184  # J139: 3312.7827 -> 4608.8184
$ws.Cells.Item(139, 11).Value = 5287.0587  # K139: 5467.5 -> 5287.0587
$ws.Cells.Item(139, 12).Value = 13826.4552  # L139: 9938.348100000001 -> 13826.4552
$ws.Cells.Item(139, 13).Value = -147.0587000000005  # M139: -327.5 -> -147.0587000000005
$ws.Cells.Item(139, 14).Value = -24106.4552  # N139: -20218.3481 -> -24106.4552

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2668.3333  # H132: 2671.5435 -> 2668.3333
$ws.Cells.Item(132, 9).Value = 1199.2  # I132: 2291 -> 1199.2
$ws.Cells.Item(132, 10).Value = 10014  # J132: 4791.7144 -> 10014
$ws.Cells.Item(132, 11).Value = 3597.6  # K132: 6873 -> 3597.6
$ws.Cells.Item(132, 12).Value = 30042  # L132: 14375.1432 -> 30042
$ws.Cells.Item(132, 13).Value = -1067.6  # M132: -4343 -> -1067.6
$ws.Cells.Item(132, 14).Value = -35102  # N132: -19435.1432 -> -35102

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 15650.777  # H22: 10181.909 -> 15650.777
$ws.Cells.Item(22, 9).Value = 950  # I22: 499.2 -> 950
$ws.Cells.Item(22, 10).Value = 17488.375  # J22: 18250.834 -> 17488.375
$ws.Cells.Item(22, 11).Value = 950  # K22: 499.2 -> 950
$ws.Cells.Item(22, 12).Value = 17488.375  # L22: 18250.834 -> 17488.375
$ws.Cells.Item(22, 13).Value = -655  # M22: -204.2 -> -655
$ws.Cells.Item(22, 14).Value = -18078.375  # N22: -18840.834 -> -18078.375

$ws.Cells.Item(27, 8).Value = 15650.777  # H27: 10181.909 -> 15650.777
$ws.Cells.Item(27, 9).Value = 950  # I27: 499.2 -> 950
$ws.Cells.Item(27, 10).Value = 17488.375  # J27: 18250.834 -> 17488.375
$ws.Cells.Item(27, 11).Value = 950  # K27: 499.2 -> 950
$ws.Cells.Item(27, 12).Value = 17488.375  # L27: 18250.834 -> 17488.375
$ws.Cells.Item(27, 13).Value = -843  # M27: -392.2 -> -843
$ws.Cells.Item(27, 14).Value = -17702.375  # N27: -18464.834 -> -17702.375

$ws.Cells.Item(46, 8).Value = 603.3333  # H46: 610 -> 603.3333
$ws.Cells.Item(46, 9).Value = 603.3333  # I46: 610 -> 603.3333
$ws.Cells.Item(46, 11).Value = 603.3333  # K46: 610 -> 603.3333
$ws.Cells.Item(46, 13).Value = -415.3333  # M46: -422 -> -415.3333

$ws.Cells.Item(68, 8).Value = 1486.3422  # H68: 1489.2051 -> 1486.3422
$ws.Cells.Item(68, 9).Value = 1441.8485  # I68: 1446.4412 -> 1441.8485
$ws.Cells.Item(68, 11).Value = 1441.8485  # K68: 1446.4412 -> 1441.8485
$ws.Cells.Item(68, 13).Value = -692.8485000000001  # M68: -697.4412 -> -692.8485000000001

$ws.Cells.Item(71, 8).Value = 1486.3422  # H71: 1489.2051 -> 1486.3422
$ws.Cells.Item(71, 9).Value = 1441.8485  # I71: 1446.4412 -> 1441.8485
$ws.Cells.Item(71, 11).Value = 7209.2425  # K71: 7232.206 -> 7209.2425
$ws.Cells.Item(71, 13).Value = -3465.2425  # M71: -3488.206 -> -3465.2425

$ws.Cells.Item(132, 8).Value = 2289.378  # H132: 2543.9155 -> 2289.378
$ws.Cells.Item(132, 9).Value = 1972.7122  # I132: 2249.875 -> 1972.7122
$ws.Cells.Item(132, 10).Value = 3595.625  # J132: 3641.6667 -> 3595.625
$ws.Cells.Item(132, 11).Value = 5918.1366  # K132: 6749.625 -> 5918.1366
$ws.Cells.Item(132, 12).Value = 10786.875  # L132: 10925.0001 -> 10786.875
$ws.Cells.Item(132, 13).Value = -3388.1366  # M132: -4219.625 -> -3388.1366
$ws.Cells.Item(132, 14).Value = -15846.875  # N132: -15985.0001 -> -15846.875

$ws.Cells.Item(136, 8).Value = 4168085.5  # H136: 5954086.5 -> 4168085.5
$ws.Cells.Item(136, 9).Value = 1313.3704  # I136: 1692 -> 1313.3704
$ws.Cells.Item(136, 10).Value = 12822150  # J136: 15153241 -> 12822150
$ws.Cells.Item(136, 11).Value = 3940.1112  # K136: 5076 -> 3940.1112
$ws.Cells.Item(136, 12).Value = 38466450  # L136: 45459723 -> 38466450
$ws.Cells.Item(136, 13).Value = -1390.1112  # M136: -2526 -> -1390.1112
$ws.Cells.Item(136, 14).Value = -38471550  # N136: -45464823 -> -38471550

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 3268949  # H132: 3206085.2 -> 3268949
$ws.Cells.Item(132, 9).Value = 879.1389  # I132: 817.64105 -> 879.1389
$ws.Cells.Item(132, 10).Value = 11112316  # J132: 12821888 -> 11112316
$ws.Cells.Item(132, 11).Value = 2637.4167  # K132: 2452.92315 -> 2637.4167
$ws.Cells.Item(132, 12).Value = 33336948  # L132: 38465664 -> 33336948
$ws.Cells.Item(132, 13).Value = -107.4167000000002  # M132: 77.07685000000038 -> -107.4167000000002
$ws.Cells.Item(132, 14).Value = -33342008  # N132: -38470724 -> -33342008

$ws.Cells.Item(136, 8).Value = 1949.115  # H136: 1904.914 -> 1949.115
$ws.Cells.Item(136, 9).Value = 1740.0299  # I136: 1732.9559 -> 1740.0299
$ws.Cells.Item(136, 10).Value = 2649.55  # J136: 2372.64 -> 2649.55
$ws.Cells.Item(136, 11).Value = 5220.0897  # K136: 5198.8677 -> 5220.0897
$ws.Cells.Item(136, 12).Value = 7948.650000000001  # L136: 7117.92 -> 7948.650000000001
$ws.Cells.Item(136, 13).Value = -2670.0897  # M136: -2648.8677 -> -2670.0897
$ws.Cells.Item(136, 14).Value = -13048.65  # N136: -12217.92 -> -13048.65
